# Auto-generated: apply scheduled market-data refresh to Sheets workbook
# Updates currentAveragePrice* / Leve*Price / LeveProfit* columns (H:N) per-row
# across all 8 job sheets, matching the upstream scheduled-runner commit.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 4214.6665   # H9: 4346.9688 -> 4214.6665
$ws.Cells.Item(9, 9).Value = 5179.269   # I9: 5180.0386 -> 5179.269
$ws.Cells.Item(9, 10).Value = 631.8570999999999   # J9: 737 -> 631.8570999999999
$ws.Cells.Item(9, 11).Value = 5179.269   # K9: 5180.0386 -> 5179.269
$ws.Cells.Item(9, 12).Value = 631.8570999999999   # L9: 737 -> 631.8570999999999
$ws.Cells.Item(9, 13).Value = -5010.269   # M9: -5011.0386 -> -5010.269
$ws.Cells.Item(9, 14).Value = -969.8570999999999   # N9: -1075 -> -969.8570999999999

$ws.Cells.Item(40, 8).Value = 8865.666999999999   # H40: 8414.462 -> 8865.666999999999
$ws.Cells.Item(40, 10).Value = 11925   # J40: 10140 -> 11925
$ws.Cells.Item(40, 12).Value = 11925   # L40: 10140 -> 11925
$ws.Cells.Item(40, 14).Value = -12275   # N40: -10490 -> -12275

$ws.Cells.Item(47, 8).Value = 49999   # H47: 8000 -> 49999
$ws.Cells.Item(47, 9).Value = 49999   # I47: 8000 -> 49999
$ws.Cells.Item(47, 11).Value = 49999   # K47: 8000 -> 49999
$ws.Cells.Item(47, 13).Value = -49027   # M47: -7028 -> -49027

$ws.Cells.Item(51, 8).Value = 8723.333000000001   # H51: 9056.666999999999 -> 8723.333000000001
$ws.Cells.Item(51, 9).Value = 8112.6665   # I51: 9170 -> 8112.6665
$ws.Cells.Item(51, 10).Value = 9334   # J51: 9000 -> 9334
$ws.Cells.Item(51, 11).Value = 8112.6665   # K51: 9170 -> 8112.6665
$ws.Cells.Item(51, 12).Value = 9334   # L51: 9000 -> 9334
$ws.Cells.Item(51, 13).Value = -7628.6665   # M51: -8686 -> -7628.6665
$ws.Cells.Item(51, 14).Value = -10302   # N51: -9968 -> -10302

$ws.Cells.Item(94, 8).Value = 1502.8889   # H94: 1503.8889 -> 1502.8889
$ws.Cells.Item(94, 9).Value = 503.85715   # I94: 505.14285 -> 503.85715
$ws.Cells.Item(94, 11).Value = 503.85715   # K94: 505.14285 -> 503.85715
$ws.Cells.Item(94, 13).Value = -52.85714999999999   # M94: -54.14285000000001 -> -52.85714999999999

$ws.Cells.Item(98, 8).Value = 2145.9524   # H98: 2146.0977 -> 2145.9524
$ws.Cells.Item(98, 9).Value = 1850.6451   # I98: 1853.9656 -> 1850.6451
$ws.Cells.Item(98, 10).Value = 2978.182   # J98: 2852.0833 -> 2978.182
$ws.Cells.Item(98, 11).Value = 1850.6451   # K98: 1853.9656 -> 1850.6451
$ws.Cells.Item(98, 12).Value = 2978.182   # L98: 2852.0833 -> 2978.182
$ws.Cells.Item(98, 13).Value = -352.6451   # M98: -355.9656 -> -352.6451
$ws.Cells.Item(98, 14).Value = -5974.182   # N98: -5848.0833 -> -5974.182

$ws.Cells.Item(108, 8).Value = 49000   # H108: 48969.332 -> 49000
$ws.Cells.Item(108, 10).Value = 49000   # J108: 48969.332 -> 49000
$ws.Cells.Item(108, 12).Value = 49000   # L108: 48969.332 -> 49000
$ws.Cells.Item(108, 14).Value = -56680   # N108: -56649.332 -> -56680

$ws.Cells.Item(110, 8).Value = 200700   # H110: 0 -> 200700
$ws.Cells.Item(110, 10).Value = 200700   # J110: 0 -> 200700
$ws.Cells.Item(110, 12).Value = 200700   # L110: 0 -> 200700
$ws.Cells.Item(110, 14).Value = -208880   # N110: None -> -208880

$ws.Cells.Item(122, 8).Value = 2145.9524   # H122: 2146.0977 -> 2145.9524
$ws.Cells.Item(122, 9).Value = 1850.6451   # I122: 1853.9656 -> 1850.6451
$ws.Cells.Item(122, 10).Value = 2978.182   # J122: 2852.0833 -> 2978.182
$ws.Cells.Item(122, 11).Value = 5551.9353   # K122: 5561.8968 -> 5551.9353
$ws.Cells.Item(122, 12).Value = 8934.545999999998   # L122: 8556.249899999999 -> 8934.545999999998
$ws.Cells.Item(122, 13).Value = -3101.9353   # M122: -3111.8968 -> -3101.9353
$ws.Cells.Item(122, 14).Value = -13834.546   # N122: -13456.2499 -> -13834.546

$ws.Cells.Item(132, 8).Value = 2192.0625   # H132: 2263.4355 -> 2192.0625
$ws.Cells.Item(132, 9).Value = 2157.9673   # I132: 2328.1272 -> 2157.9673
$ws.Cells.Item(132, 10).Value = 2885.3333   # J132: 1755.1428 -> 2885.3333
$ws.Cells.Item(132, 11).Value = 6473.901899999999   # K132: 6984.3816 -> 6473.901899999999
$ws.Cells.Item(132, 12).Value = 8655.999899999999   # L132: 5265.428400000001 -> 8655.999899999999
$ws.Cells.Item(132, 13).Value = -3943.901899999999   # M132: -4454.3816 -> -3943.901899999999
$ws.Cells.Item(132, 14).Value = -13715.9999   # N132: -10325.4284 -> -13715.9999

$ws.Cells.Item(138, 8).Value = 6723.9707   # H138: 6100.2925 -> 6723.9707
$ws.Cells.Item(138, 9).Value = 3421.842   # I138: 3295.6 -> 3421.842
$ws.Cells.Item(138, 10).Value = 10906.667   # J138: 8771.429 -> 10906.667
$ws.Cells.Item(138, 11).Value = 10265.526   # K138: 9886.799999999999 -> 10265.526
$ws.Cells.Item(138, 12).Value = 32720.001   # L138: 26314.287 -> 32720.001
$ws.Cells.Item(138, 13).Value = -5125.526   # M138: -4746.799999999999 -> -5125.526
$ws.Cells.Item(138, 14).Value = -43000.001   # N138: -36594.287 -> -43000.001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3714.772   # H32: 3411.0984 -> 3714.772
$ws.Cells.Item(32, 9).Value = 3173.9644   # I32: 2901.2834 -> 3173.9644
$ws.Cells.Item(32, 11).Value = 3173.9644   # K32: 2901.2834 -> 3173.9644
$ws.Cells.Item(32, 13).Value = -2886.9644   # M32: -2614.2834 -> -2886.9644

$ws.Cells.Item(45, 8).Value = 5927.8125   # H45: 6308 -> 5927.8125
$ws.Cells.Item(45, 9).Value = 4843.143   # I45: 5612.8335 -> 4843.143
$ws.Cells.Item(45, 11).Value = 4843.143   # K45: 5612.8335 -> 4843.143
$ws.Cells.Item(45, 13).Value = -4466.143   # M45: -5235.8335 -> -4466.143

$ws.Cells.Item(61, 8).Value = 3217.6853   # H61: 3188.3167 -> 3217.6853
$ws.Cells.Item(61, 9).Value = 2227.3953   # I61: 2202.932 -> 2227.3953
$ws.Cells.Item(61, 10).Value = 7088.8184   # J61: 5898.125 -> 7088.8184
$ws.Cells.Item(61, 11).Value = 2227.3953   # K61: 2202.932 -> 2227.3953
$ws.Cells.Item(61, 12).Value = 7088.8184   # L61: 5898.125 -> 7088.8184
$ws.Cells.Item(61, 13).Value = -2015.3953   # M61: -1990.932 -> -2015.3953
$ws.Cells.Item(61, 14).Value = -7512.8184   # N61: -6322.125 -> -7512.8184

$ws.Cells.Item(97, 8).Value = 1889.04   # H97: 1800.4642 -> 1889.04
$ws.Cells.Item(97, 9).Value = 1882.6818   # I97: 1886.7727 -> 1882.6818
$ws.Cells.Item(97, 10).Value = 1935.6666   # J97: 1484 -> 1935.6666
$ws.Cells.Item(97, 11).Value = 1882.6818   # K97: 1886.7727 -> 1882.6818
$ws.Cells.Item(97, 12).Value = 1935.6666   # L97: 1484 -> 1935.6666
$ws.Cells.Item(97, 13).Value = -1386.6818   # M97: -1390.7727 -> -1386.6818
$ws.Cells.Item(97, 14).Value = -2927.6666   # N97: -2476 -> -2927.6666

$ws.Cells.Item(110, 8).Value = 313480.25   # H110: 455724.8 -> 313480.25
$ws.Cells.Item(110, 9).Value = 357941   # I110: 500996.3 -> 357941
$ws.Cells.Item(110, 10).Value = 2255   # J110: 3010 -> 2255
$ws.Cells.Item(110, 11).Value = 357941   # K110: 500996.3 -> 357941
$ws.Cells.Item(110, 12).Value = 2255   # L110: 3010 -> 2255
$ws.Cells.Item(110, 13).Value = -355896   # M110: -498951.3 -> -355896
$ws.Cells.Item(110, 14).Value = -6345   # N110: -7100 -> -6345

$ws.Cells.Item(122, 8).Value = 2401.4   # H122: 2829.0476 -> 2401.4
$ws.Cells.Item(122, 9).Value = 1261.7428   # I122: 1541.4642 -> 1261.7428
$ws.Cells.Item(122, 10).Value = 5060.6   # J122: 5404.2144 -> 5060.6
$ws.Cells.Item(122, 11).Value = 3785.2284   # K122: 4624.392599999999 -> 3785.2284
$ws.Cells.Item(122, 12).Value = 15181.8   # L122: 16212.6432 -> 15181.8
$ws.Cells.Item(122, 13).Value = -1335.2284   # M122: -2174.392599999999 -> -1335.2284
$ws.Cells.Item(122, 14).Value = -20081.8   # N122: -21112.6432 -> -20081.8

$ws.Cells.Item(136, 8).Value = 3217.6853   # H136: 3188.3167 -> 3217.6853
$ws.Cells.Item(136, 9).Value = 2227.3953   # I136: 2202.932 -> 2227.3953
$ws.Cells.Item(136, 10).Value = 7088.8184   # J136: 5898.125 -> 7088.8184
$ws.Cells.Item(136, 11).Value = 6682.1859   # K136: 6608.795999999999 -> 6682.1859
$ws.Cells.Item(136, 12).Value = 21266.4552   # L136: 17694.375 -> 21266.4552
$ws.Cells.Item(136, 13).Value = -4132.1859   # M136: -4058.795999999999 -> -4132.1859
$ws.Cells.Item(136, 14).Value = -26366.4552   # N136: -22794.375 -> -26366.4552

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 3229.8696   # H107: 3555.8333 -> 3229.8696
$ws.Cells.Item(107, 9).Value = 3063.7   # I107: 3467 -> 3063.7
$ws.Cells.Item(107, 10).Value = 4337.6665   # J107: 4000 -> 4337.6665
$ws.Cells.Item(107, 11).Value = 3063.7   # K107: 3467 -> 3063.7
$ws.Cells.Item(107, 12).Value = 4337.6665   # L107: 4000 -> 4337.6665
$ws.Cells.Item(107, 13).Value = -1143.7   # M107: -1547 -> -1143.7
$ws.Cells.Item(107, 14).Value = -8177.6665   # N107: -7840 -> -8177.6665

$ws.Cells.Item(134, 8).Value = 22809.203   # H134: 22823.314 -> 22809.203
$ws.Cells.Item(134, 9).Value = 3064.743   # I134: 3119.2354 -> 3064.743
$ws.Cells.Item(134, 10).Value = 59180.58   # J134: 56320.25 -> 59180.58
$ws.Cells.Item(134, 11).Value = 9194.228999999999   # K134: 9357.706200000001 -> 9194.228999999999
$ws.Cells.Item(134, 12).Value = 177541.74   # L134: 168960.75 -> 177541.74
$ws.Cells.Item(134, 13).Value = -6659.228999999999   # M134: -6822.706200000001 -> -6659.228999999999
$ws.Cells.Item(134, 14).Value = -182611.74   # N134: -174030.75 -> -182611.74

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 4617.59   # H16: 4621.15 -> 4617.59
$ws.Cells.Item(16, 9).Value = 2412.5652   # I16: 2554.4783 -> 2412.5652
$ws.Cells.Item(16, 10).Value = 7787.3125   # J16: 7417.2354 -> 7787.3125
$ws.Cells.Item(16, 11).Value = 2412.5652   # K16: 2554.4783 -> 2412.5652
$ws.Cells.Item(16, 12).Value = 7787.3125   # L16: 7417.2354 -> 7787.3125
$ws.Cells.Item(16, 13).Value = -2125.5652   # M16: -2267.4783 -> -2125.5652
$ws.Cells.Item(16, 14).Value = -8361.3125   # N16: -7991.2354 -> -8361.3125

$ws.Cells.Item(22, 8).Value = 449.33334   # H22: 479.4 -> 449.33334
$ws.Cells.Item(22, 9).Value = 339.2   # I22: 349.25 -> 339.2
$ws.Cells.Item(22, 11).Value = 339.2   # K22: 349.25 -> 339.2
$ws.Cells.Item(22, 13).Value = 10.80000000000001   # M22: 0.75 -> 10.80000000000001

$ws.Cells.Item(31, 8).Value = 361219.6   # H31: 440429.7 -> 361219.6
$ws.Cells.Item(31, 9).Value = 626804.9   # I31: 1113298.6 -> 626804.9
$ws.Cells.Item(31, 10).Value = 7105.9165   # J31: 7871.0713 -> 7105.9165
$ws.Cells.Item(31, 11).Value = 626804.9   # K31: 1113298.6 -> 626804.9
$ws.Cells.Item(31, 12).Value = 7105.9165   # L31: 7871.0713 -> 7105.9165
$ws.Cells.Item(31, 13).Value = -626509.9   # M31: -1113003.6 -> -626509.9
$ws.Cells.Item(31, 14).Value = -7695.9165   # N31: -8461.0713 -> -7695.9165

$ws.Cells.Item(34, 8).Value = 361219.6   # H34: 440429.7 -> 361219.6
$ws.Cells.Item(34, 9).Value = 626804.9   # I34: 1113298.6 -> 626804.9
$ws.Cells.Item(34, 10).Value = 7105.9165   # J34: 7871.0713 -> 7105.9165
$ws.Cells.Item(34, 11).Value = 626804.9   # K34: 1113298.6 -> 626804.9
$ws.Cells.Item(34, 12).Value = 7105.9165   # L34: 7871.0713 -> 7105.9165
$ws.Cells.Item(34, 13).Value = -626602.9   # M34: -1113096.6 -> -626602.9
$ws.Cells.Item(34, 14).Value = -7509.9165   # N34: -8275.0713 -> -7509.9165

$ws.Cells.Item(58, 8).Value = 189349.22   # H58: 189312.22 -> 189349.22
$ws.Cells.Item(58, 9).Value = 387173.28   # I58: 419174.9 -> 387173.28
$ws.Cells.Item(58, 10).Value = 5655.4644   # J58: 5422.067 -> 5655.4644
$ws.Cells.Item(58, 11).Value = 387173.28   # K58: 419174.9 -> 387173.28
$ws.Cells.Item(58, 12).Value = 5655.4644   # L58: 5422.067 -> 5655.4644
$ws.Cells.Item(58, 13).Value = -386970.28   # M58: -418971.9 -> -386970.28
$ws.Cells.Item(58, 14).Value = -6061.4644   # N58: -5828.067 -> -6061.4644

$ws.Cells.Item(62, 8).Value = 4182.1665   # H62: 4364.2 -> 4182.1665
$ws.Cells.Item(62, 9).Value = 3370.3333   # I62: 3393.4 -> 3370.3333
$ws.Cells.Item(62, 10).Value = 4994   # J62: 5335 -> 4994
$ws.Cells.Item(62, 11).Value = 3370.3333   # K62: 3393.4 -> 3370.3333
$ws.Cells.Item(62, 12).Value = 4994   # L62: 5335 -> 4994
$ws.Cells.Item(62, 13).Value = -2746.3333   # M62: -2769.4 -> -2746.3333
$ws.Cells.Item(62, 14).Value = -6242   # N62: -6583 -> -6242

$ws.Cells.Item(65, 8).Value = 4182.1665   # H65: 4364.2 -> 4182.1665
$ws.Cells.Item(65, 9).Value = 3370.3333   # I65: 3393.4 -> 3370.3333
$ws.Cells.Item(65, 10).Value = 4994   # J65: 5335 -> 4994
$ws.Cells.Item(65, 11).Value = 16851.6665   # K65: 16967 -> 16851.6665
$ws.Cells.Item(65, 12).Value = 24970   # L65: 26675 -> 24970
$ws.Cells.Item(65, 13).Value = -13731.6665   # M65: -13847 -> -13731.6665
$ws.Cells.Item(65, 14).Value = -31210   # N65: -32915 -> -31210

$ws.Cells.Item(76, 8).Value = 6666.3335   # H76: 7500 -> 6666.3335
$ws.Cells.Item(76, 9).Value = 6666.3335   # I76: 7500 -> 6666.3335
$ws.Cells.Item(76, 11).Value = 6666.3335   # K76: 7500 -> 6666.3335
$ws.Cells.Item(76, 13).Value = -6351.3335   # M76: -7185 -> -6351.3335

$ws.Cells.Item(79, 8).Value = 6666.3335   # H79: 7500 -> 6666.3335
$ws.Cells.Item(79, 9).Value = 6666.3335   # I79: 7500 -> 6666.3335
$ws.Cells.Item(79, 11).Value = 6666.3335   # K79: 7500 -> 6666.3335
$ws.Cells.Item(79, 13).Value = -5574.3335   # M79: -6408 -> -5574.3335

$ws.Cells.Item(94, 8).Value = 1494.1765   # H94: 1559.875 -> 1494.1765
$ws.Cells.Item(94, 10).Value = 1342.8182   # J94: 1432.8 -> 1342.8182
$ws.Cells.Item(94, 12).Value = 1342.8182   # L94: 1432.8 -> 1342.8182
$ws.Cells.Item(94, 14).Value = -2244.8182   # N94: -2334.8 -> -2244.8182

$ws.Cells.Item(105, 8).Value = 1692.0769   # H105: 1580.5714 -> 1692.0769
$ws.Cells.Item(105, 9).Value = 941.6667   # I105: 876.2857 -> 941.6667
$ws.Cells.Item(105, 10).Value = 2335.2856   # J105: 2284.8572 -> 2335.2856
$ws.Cells.Item(105, 11).Value = 941.6667   # K105: 876.2857 -> 941.6667
$ws.Cells.Item(105, 12).Value = 2335.2856   # L105: 2284.8572 -> 2335.2856
$ws.Cells.Item(105, 13).Value = 805.3333   # M105: 870.7143 -> 805.3333
$ws.Cells.Item(105, 14).Value = -5829.2856   # N105: -5778.8572 -> -5829.2856

$ws.Cells.Item(113, 8).Value = 4617.59   # H113: 4621.15 -> 4617.59
$ws.Cells.Item(113, 9).Value = 2412.5652   # I113: 2554.4783 -> 2412.5652
$ws.Cells.Item(113, 10).Value = 7787.3125   # J113: 7417.2354 -> 7787.3125
$ws.Cells.Item(113, 11).Value = 2412.5652   # K113: 2554.4783 -> 2412.5652
$ws.Cells.Item(113, 12).Value = 7787.3125   # L113: 7417.2354 -> 7787.3125
$ws.Cells.Item(113, 13).Value = -242.5652   # M113: -384.4783000000002 -> -242.5652
$ws.Cells.Item(113, 14).Value = -12127.3125   # N113: -11757.2354 -> -12127.3125

$ws.Cells.Item(132, 8).Value = 3656.7104   # H132: 3583.625 -> 3656.7104
$ws.Cells.Item(132, 9).Value = 2263   # I132: 2258.1428 -> 2263
$ws.Cells.Item(132, 11).Value = 6789   # K132: 6774.428400000001 -> 6789
$ws.Cells.Item(132, 13).Value = -4259   # M132: -4244.428400000001 -> -4259

$ws.Cells.Item(136, 8).Value = 189349.22   # H136: 189312.22 -> 189349.22
$ws.Cells.Item(136, 9).Value = 387173.28   # I136: 419174.9 -> 387173.28
$ws.Cells.Item(136, 10).Value = 5655.4644   # J136: 5422.067 -> 5655.4644
$ws.Cells.Item(136, 11).Value = 1161519.84   # K136: 1257524.7 -> 1161519.84
$ws.Cells.Item(136, 12).Value = 16966.3932   # L136: 16266.201 -> 16966.3932
$ws.Cells.Item(136, 13).Value = -1158969.84   # M136: -1254974.7 -> -1158969.84
$ws.Cells.Item(136, 14).Value = -22066.3932   # N136: -21366.201 -> -22066.3932

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(14, 8).Value = 339.2   # H14: 1217.4546 -> 339.2
$ws.Cells.Item(14, 9).Value = 339.2   # I14: 1217.4546 -> 339.2
$ws.Cells.Item(14, 11).Value = 1017.6   # K14: 3652.3638 -> 1017.6
$ws.Cells.Item(14, 13).Value = -844.5999999999999   # M14: -3479.3638 -> -844.5999999999999

$ws.Cells.Item(20, 8).Value = 285   # H20: 550 -> 285
$ws.Cells.Item(20, 9).Value = 99   # I20: 100 -> 99
$ws.Cells.Item(20, 10).Value = 750   # J20: 1000 -> 750
$ws.Cells.Item(20, 11).Value = 297   # K20: 300 -> 297
$ws.Cells.Item(20, 12).Value = 2250   # L20: 3000 -> 2250
$ws.Cells.Item(20, 13).Value = -70   # M20: -73 -> -70
$ws.Cells.Item(20, 14).Value = -2704   # N20: -3454 -> -2704

$ws.Cells.Item(131, 8).Value = 3763.9678   # H131: 3713.1936 -> 3763.9678
$ws.Cells.Item(131, 9).Value = 1271.8   # I131: 1258.8823 -> 1271.8
$ws.Cells.Item(131, 10).Value = 6100.375   # J131: 6693.4287 -> 6100.375
$ws.Cells.Item(131, 11).Value = 3815.4   # K131: 3776.6469 -> 3815.4
$ws.Cells.Item(131, 12).Value = 18301.125   # L131: 20080.2861 -> 18301.125
$ws.Cells.Item(131, 13).Value = 1224.6   # M131: 1263.3531 -> 1224.6
$ws.Cells.Item(131, 14).Value = -28381.125   # N131: -30160.2861 -> -28381.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 131556.86   # H132: 134725.36 -> 131556.86
$ws.Cells.Item(132, 9).Value = 165754.02   # I132: 174307.2 -> 165754.02
$ws.Cells.Item(132, 10).Value = 44639.082   # J132: 42895.48 -> 44639.082
$ws.Cells.Item(132, 11).Value = 497262.0599999999   # K132: 522921.6 -> 497262.0599999999
$ws.Cells.Item(132, 12).Value = 133917.246   # L132: 128686.44 -> 133917.246
$ws.Cells.Item(132, 13).Value = -494732.0599999999   # M132: -520391.6 -> -494732.0599999999
$ws.Cells.Item(132, 14).Value = -138977.246   # N132: -133746.44 -> -138977.246

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 3712.125   # H46: 3649.625 -> 3712.125
$ws.Cells.Item(46, 9).Value = 2973.4   # I46: 2880.6875 -> 2973.4
$ws.Cells.Item(46, 10).Value = 4943.3335   # J46: 5187.5 -> 4943.3335
$ws.Cells.Item(46, 11).Value = 2973.4   # K46: 2880.6875 -> 2973.4
$ws.Cells.Item(46, 12).Value = 4943.3335   # L46: 5187.5 -> 4943.3335
$ws.Cells.Item(46, 13).Value = -2785.4   # M46: -2692.6875 -> -2785.4
$ws.Cells.Item(46, 14).Value = -5319.3335   # N46: -5563.5 -> -5319.3335

$ws.Cells.Item(132, 8).Value = 3109.6042   # H132: 3667.1667 -> 3109.6042
$ws.Cells.Item(132, 9).Value = 2657.6365   # I132: 3378.0881 -> 2657.6365
$ws.Cells.Item(132, 10).Value = 4103.933   # J132: 4369.2144 -> 4103.933
$ws.Cells.Item(132, 11).Value = 7972.9095   # K132: 10134.2643 -> 7972.9095
$ws.Cells.Item(132, 12).Value = 12311.799   # L132: 13107.6432 -> 12311.799
$ws.Cells.Item(132, 13).Value = -5442.9095   # M132: -7604.264299999999 -> -5442.9095
$ws.Cells.Item(132, 14).Value = -17371.799   # N132: -18167.6432 -> -17371.799

$ws.Cells.Item(136, 8).Value = 3837.1924   # H136: 3311 -> 3837.1924
$ws.Cells.Item(136, 9).Value = 3124.0588   # I136: 2706.125 -> 3124.0588
$ws.Cells.Item(136, 10).Value = 5184.222   # J136: 4630.727 -> 5184.222
$ws.Cells.Item(136, 11).Value = 9372.1764   # K136: 8118.375 -> 9372.1764
$ws.Cells.Item(136, 12).Value = 15552.666   # L136: 13892.181 -> 15552.666
$ws.Cells.Item(136, 13).Value = -6822.1764   # M136: -5568.375 -> -6822.1764
$ws.Cells.Item(136, 14).Value = -20652.666   # N136: -18992.181 -> -20652.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(41, 8).Value = 9275.888999999999   # H41: 9336.875 -> 9275.888999999999
$ws.Cells.Item(41, 10).Value = 8988.125   # J41: 9016.714 -> 8988.125
$ws.Cells.Item(41, 12).Value = 8988.125   # L41: 9016.714 -> 8988.125
$ws.Cells.Item(41, 14).Value = -9768.125   # N41: -9796.714 -> -9768.125
